$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2-5) get cyclically rotated: row2's data moves to row5,
# and rows 3,4,5 each shift up by one row (3->2, 4->3, 5->4), for
# columns A, B, E, F, G, H, Q, R.

# Capture original values before overwriting anything.
$cols = @("A","B","E","F","G","H","Q","R")
$orig = @{}
foreach ($r in 2..5) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row r gets the data that used to be in row (r+1), except row 5
# which gets what used to be in row 2.
foreach ($r in 2..5) {
    $srcRow = $r + 1
    if ($srcRow -gt 5) { $srcRow = 2 }
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig[$srcRow][$c]
    }
}
